$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the results table (rows 2-8, columns B:J) with the new benchmark numbers ---
$ws.Range("B2").Value = 6900
$ws.Range("C2").Value = 13720
$ws.Range("D2").Value = 28380
$ws.Range("E2").Value = 69960
$ws.Range("F2").Value = 154600
$ws.Range("G2").Value = 432150
$ws.Range("H2").Value = 1347540
$ws.Range("I2").Value = 4532630
$ws.Range("J2").Value = 16446350

$ws.Range("B3").Value = 7470
$ws.Range("C3").Value = 13820
$ws.Range("D3").Value = 28740
$ws.Range("E3").Value = 62330
$ws.Range("F3").Value = 149580
$ws.Range("G3").Value = 398580
$ws.Range("H3").Value = 1185890
$ws.Range("I3").Value = 3916180
$ws.Range("J3").Value = 14205750

$ws.Range("B4").Value = 5900
$ws.Range("C4").Value = 9900
$ws.Range("D4").Value = 17640
$ws.Range("E4").Value = 32910
$ws.Range("F4").Value = 61960
$ws.Range("G4").Value = 123790
$ws.Range("H4").Value = 233340
$ws.Range("I4").Value = 463230
$ws.Range("J4").Value = 916310

$ws.Range("B5").Value = 8160
$ws.Range("C5").Value = 18020
$ws.Range("D5").Value = 39220
$ws.Range("E5").Value = 84860
$ws.Range("F5").Value = 180460
$ws.Range("G5").Value = 378480
$ws.Range("H5").Value = 757480
$ws.Range("I5").Value = 1508130
$ws.Range("J5").Value = 3012890

$ws.Range("B6").Value = 75050
$ws.Range("C6").Value = 291030
$ws.Range("D6").Value = 1165580
$ws.Range("E6").Value = 4668220
$ws.Range("F6").Value = 18541010
$ws.Range("G6").Value = 74027310
$ws.Range("H6").Value = 294371060
$ws.Range("I6").Value = 1183581160
$ws.Range("J6").Value = 4725895550

$ws.Range("B7").Value = 5480
$ws.Range("C7").Value = 9780
$ws.Range("D7").Value = 17600
$ws.Range("E7").Value = 33050
$ws.Range("F7").Value = 63930
$ws.Range("G7").Value = 123440
$ws.Range("H7").Value = 241910
$ws.Range("I7").Value = 477680
$ws.Range("J7").Value = 944370

$ws.Range("B8").Value = 5620
$ws.Range("C8").Value = 10970
$ws.Range("D8").Value = 21310
$ws.Range("E8").Value = 42120
$ws.Range("F8").Value = 83760
$ws.Range("G8").Value = 164110
$ws.Range("H8").Value = 339750
$ws.Range("I8").Value = 691640
$ws.Range("J8").Value = 1415520

# --- Rename the two chart titles ---
$co1 = $ws.ChartObjects(1)
$co2 = $ws.ChartObjects(2)
$co1.Chart.ChartTitle.Text = "Average-Cases"
$co2.Chart.ChartTitle.Text = "Average-Cases ohne InsertionSort"

# --- Reposition / resize the two charts ---
$co1.Top = 156.37496062992125
$co1.Left = 16.12488188976378
$co1.Width = 672.3038290477363
$co1.Height = 340.87503937007875

$co2.Top = 156.37488188976377
$co2.Left = 693.6786321973425
$co2.Width = 537.1875787401575
$co2.Height = 338.62511811023626

# --- Update the active cell selection ---
$ws.Range("T20").Select() | Out-Null
